$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1 (paragraph "There must be something here...Authors"):
# merge the two runs ("...about the " + "Authors") into a single run and
# drop the gramStart/gramEnd proofErr markers that bracketed "Authors".
# A plain text/range edit leaves the proofErr markers orphaned in place, so
# instead we delete the whole paragraph (including its mark) and recreate it
# fresh with a single run - that naturally drops the stale proofErr markers.
# ---------------------------------------------------------------------------
$mergedSentence = "There must be something here- ----- I couldn" + [char]0x2019 + "t understand what kind of thing need to be included about the Authors"

$paraIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "There must be something here*Authors*") {
        $paraIndex = $i
        break
    }
}

if ($paraIndex -gt 0) {
    $prevPara = $d.Paragraphs($paraIndex - 1)
    $target = $d.Paragraphs($paraIndex)
    $full = $d.Range($target.Range.Start, $target.Range.End)
    $full.Delete()
    $prevPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs($paraIndex)
    $newPara.Range.InsertBefore($mergedSentence)
}

# ---------------------------------------------------------------------------
# Change 2 (Abstract, first sentence of the body):
# "This paper " + "summarizes" + " the comparison of programming exercises
#  on C# providing by "
#    -> "In t" + "his paper " + "we present a comparative study" +
#       " of programming exercises on C# providing by "
#
# Word (and this engine) silently re-merges any adjacent runs that end up
# with identical final formatting, which would collapse our careful run
# split back into one big run. Toggling a character property (Bold) on/off
# around every Range.Text write keeps that edited span from being folded
# back into its neighbours, letting us reproduce the run boundaries from
# the diff exactly. (NOTE: keep every step inline - routing a Range through
# a PowerShell function here does not keep the live COM reference, so
# .Start/.End stop tracking edits made inside the function.)
# ---------------------------------------------------------------------------
$part1 = "In t"
$part2 = "his paper "
$part3 = "we present a comparative study"

$rA = $d.Content
$rA.Find.Execute("This paper summarizes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startA = $rA.Start

$rA.Font.Bold = $true
$rA.Text = ($part1 + $part2 + $part3)
$rA.Font.Bold = $false

$split1 = $startA + $part1.Length
$rB = $d.Range($split1, $rA.End)
$rB.Font.Bold = $true
$rB.Text = ($part2 + $part3)
$rB.Font.Bold = $false

$split2 = $split1 + $part2.Length
$rC = $d.Range($split2, $rB.End)
$rC.Font.Bold = $true
$rC.Text = $part3
$rC.Font.Bold = $false

$rD = $d.Content
$rD.Find.Execute(" the comparison of programming exercises on C# providing by ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rD.Font.Bold = $true
$rD.Text = " of programming exercises on C# providing by "
$rD.Font.Bold = $false

# ---------------------------------------------------------------------------
# Change 3 (Abstract, end of the sentence listing the exercise chapters):
# "eritance, ... Generics Types and Class. The possible conclusion ...
#  and limitation"
#    -> "eritance, ... Generics Types and Class. " +
#       "These programming exercises are compared under the characteristics
#        of analysis, design, implementation and testing. " +
#       "The possible conclusion ... and limitation"
# ---------------------------------------------------------------------------
$seg1 = "eritance, and lastly Generics Types and Class. "
$seg2 = "These programming exercises are compared under the characteristics of analysis, design, implementation and testing. "
$seg3 = "The possible conclusion would be forwarded based on the outcome of each programming exercise, and limitation"

$rE = $d.Content
$rE.Find.Execute(
    "eritance, and lastly Generics Types and Class. The possible conclusion would be forwarded based on the outcome of each programming exercise, and limitation",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startE = $rE.Start

$rE.Font.Bold = $true
$rE.Text = ($seg1 + $seg2 + $seg3)
$rE.Font.Bold = $false

$splitE1 = $startE + $seg1.Length
$rF = $d.Range($splitE1, $rE.End)
$rF.Font.Bold = $true
$rF.Text = ($seg2 + $seg3)
$rF.Font.Bold = $false

$splitE2 = $splitE1 + $seg2.Length
$rG = $d.Range($splitE2, $rF.End)
$rG.Font.Bold = $true
$rG.Text = $seg3
$rG.Font.Bold = $false
